$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: JMISMDL (rows 434-442) ---
# Copy each source row individually (row-by-row) so Excel does not
# split the merged-cell border into separate top/middle/bottom segments.
$ws.Range("A425:G425").Copy($ws.Range("A434:G434"))
$ws.Range("A426:G426").Copy($ws.Range("A435:G435"))
$ws.Range("A427:G427").Copy($ws.Range("A436:G436"))
$ws.Range("A428:G428").Copy($ws.Range("A437:G437"))
$ws.Range("A429:G429").Copy($ws.Range("A438:G438"))
$ws.Range("A430:G430").Copy($ws.Range("A439:G439"))
$ws.Range("A431:G431").Copy($ws.Range("A440:G440"))
$ws.Range("A432:G432").Copy($ws.Range("A441:G441"))
$ws.Range("A433:G433").Copy($ws.Range("A442:G442"))
$ws.Range("A434:A442").Merge()
# Re-apply a plain thin box border per-cell so the merge operation
# reuses the single shared "box" style instead of minting top/mid/bottom variants.
$ws.Cells.Item(434, 1).Borders.LineStyle = 1
$ws.Cells.Item(435, 1).Borders.LineStyle = 1
$ws.Cells.Item(436, 1).Borders.LineStyle = 1
$ws.Cells.Item(437, 1).Borders.LineStyle = 1
$ws.Cells.Item(438, 1).Borders.LineStyle = 1
$ws.Cells.Item(439, 1).Borders.LineStyle = 1
$ws.Cells.Item(440, 1).Borders.LineStyle = 1
$ws.Cells.Item(441, 1).Borders.LineStyle = 1
$ws.Cells.Item(442, 1).Borders.LineStyle = 1

# --- Block 2: BPML (rows 443-451) ---
$ws.Range("A425:G425").Copy($ws.Range("A443:G443"))
$ws.Range("A426:G426").Copy($ws.Range("A444:G444"))
$ws.Range("A427:G427").Copy($ws.Range("A445:G445"))
$ws.Range("A428:G428").Copy($ws.Range("A446:G446"))
$ws.Range("A429:G429").Copy($ws.Range("A447:G447"))
$ws.Range("A430:G430").Copy($ws.Range("A448:G448"))
$ws.Range("A431:G431").Copy($ws.Range("A449:G449"))
$ws.Range("A432:G432").Copy($ws.Range("A450:G450"))
$ws.Range("A433:G433").Copy($ws.Range("A451:G451"))
$ws.Range("A443:A451").Merge()
$ws.Cells.Item(443, 1).Borders.LineStyle = 1
$ws.Cells.Item(444, 1).Borders.LineStyle = 1
$ws.Cells.Item(445, 1).Borders.LineStyle = 1
$ws.Cells.Item(446, 1).Borders.LineStyle = 1
$ws.Cells.Item(447, 1).Borders.LineStyle = 1
$ws.Cells.Item(448, 1).Borders.LineStyle = 1
$ws.Cells.Item(449, 1).Borders.LineStyle = 1
$ws.Cells.Item(450, 1).Borders.LineStyle = 1
$ws.Cells.Item(451, 1).Borders.LineStyle = 1

# --- Write the actual data values ---
$ws.Range("A434").Value = "JMISMDL"
$ws.Range("B434").Value = "current ratio"
$ws.Range("C434").Value = 3.358752166377816
$ws.Range("D434").Value = 5.52
$ws.Range("E434").Value = 3.553846153846154
$ws.Range("F434").Value = 5.1353591160221
$ws.Range("G434").Value = 1.943970767356882

$ws.Range("A435").ClearContents()
$ws.Range("B435").Value = "cash ratio"
$ws.Range("C435").Value = 0.6065857885615251
$ws.Range("D435").Value = 0.4166666666666667
$ws.Range("E435").Value = 0.05576923076923077
$ws.Range("F435").Value = 0.9364640883977901
$ws.Range("G435").Value = 0.01339829476248477

$ws.Range("A436").ClearContents()
$ws.Range("B436").Value = "Total Asset turnover"
$ws.Range("C436").Value = 0.7128654970760234
$ws.Range("D436").Value = 0.6397176772537696
$ws.Range("E436").Value = 0.5945864661654136
$ws.Range("F436").Value = 0.5850608752879236
$ws.Range("G436").Value = 0.5078627591136526

$ws.Range("A437").ClearContents()
$ws.Range("B437").Value = "Account receivable turnover"
$ws.Range("C437").Value = 4.670498084291188
$ws.Range("D437").Value = 3.140157480314961
$ws.Range("E437").Value = 2.254275940706956
$ws.Range("F437").Value = 3.39961759082218
$ws.Range("G437").Value = 2.722222222222222

$ws.Range("A438").ClearContents()
$ws.Range("B438").Value = "Debt Ratio"
$ws.Range("C438").Value = 0.1970760233918129
$ws.Range("D438").Value = 0.1299326275264678
$ws.Range("E438").Value = 0.1912781954887218
$ws.Range("F438").Value = 0.1523527476143468
$ws.Range("G438").Value = 0.7197998570407433

$ws.Range("A439").ClearContents()
$ws.Range("B439").Value = "Equity Ratio"
$ws.Range("C439").Value = 0.7994152046783626
$ws.Range("D439").Value = 0.8658966955405839
$ws.Range("E439").Value = 0.8051127819548872
$ws.Range("F439").Value = 0.8433695294504772
$ws.Range("G439").Value = 0.2759113652609007

$ws.Range("A440").ClearContents()
$ws.Range("B440").Value = "Profit margin ratio"
$ws.Range("C440").Value = 0.04142739950779328
$ws.Range("D440").Value = 0.04463390170511534
$ws.Range("E440").Value = 0.04855842185128983
$ws.Range("F440").Value = 0.03768278965129359
$ws.Range("G440").Value = 0.05277973258268825

$ws.Range("A441").ClearContents()
$ws.Range("B441").Value = "ROE"
$ws.Range("C441").Value = 0.03694220921726409
$ws.Range("D441").Value = 0.03297517599110782
$ws.Range("E441").Value = 0.03586103847590587
$ws.Range("F441").Value = 0.02614124073351541
$ws.Range("G441").Value = 0.09715025906735751

$ws.Range("A442").ClearContents()
$ws.Range("B442").Value = "ROA"
$ws.Range("C442").Value = 0.02953216374269006
$ws.Range("D442").Value = 0.02855309592556946
$ws.Range("E442").Value = 0.02887218045112782
$ws.Range("F442").Value = 0.02204672589667654
$ws.Range("G442").Value = 0.0268048606147248

$ws.Range("A443").Value = "BPML"
$ws.Range("B443").Value = "current ratio"
$ws.Range("C443").Value = 1.471151216844408
$ws.Range("D443").Value = 1.661256289563383
$ws.Range("E443").Value = 1.455940315315315
$ws.Range("F443").Value = 1.155278510349237
$ws.Range("G443").Value = 1.097371714643304

$ws.Range("A444").ClearContents()
$ws.Range("B444").Value = "cash ratio"
$ws.Range("C444").Value = 0.0880503144654088
$ws.Range("D444").Value = 0.09024509008277877
$ws.Range("E444").Value = 0.0736204954954955
$ws.Range("F444").Value = 0.04850955006792936
$ws.Range("G444").Value = 0.1591155611180642

$ws.Range("A445").ClearContents()
$ws.Range("B445").Value = "Total Asset turnover"
$ws.Range("C445").Value = 0.3224307109657428
$ws.Range("D445").Value = 0.3230105431643411
$ws.Range("E445").Value = 0.3101945977345338
$ws.Range("F445").Value = 0.3656993308156099
$ws.Range("G445").Value = 0.4317692123533637

$ws.Range("A446").ClearContents()
$ws.Range("B446").Value = "Account receivable turnover"
$ws.Range("C446").Value = 11.06102362204724
$ws.Range("D446").Value = 12.08723958333333
$ws.Range("E446").Value = 6.076813655761025
$ws.Range("F446").Value = 2.988359501894965
$ws.Range("G446").Value = 4.344038538739462

$ws.Range("A447").ClearContents()
$ws.Range("B447").Value = "Debt Ratio"
$ws.Range("C447").Value = 0.6250071727778734
$ws.Range("D447").Value = 0.7177354814015797
$ws.Range("E447").Value = 0.7194670345628812
$ws.Range("F447").Value = 0.7522030080169615
$ws.Range("G447").Value = 0.6987870082196154

$ws.Range("A448").ClearContents()
$ws.Range("B448").Value = "Equity Ratio"
$ws.Range("C448").Value = 0.3749928272221266
$ws.Range("D448").Value = 0.2822645185984203
$ws.Range("E448").Value = 0.2805329654371188
$ws.Range("F448").Value = 0.2477969919830385
$ws.Range("G448").Value = 0.301252892825792

$ws.Range("A449").ClearContents()
$ws.Range("B449").Value = "Profit margin ratio"
$ws.Range("C449").Value = 0.04520377291332978
$ws.Range("D449").Value = 0.0442744802326834
$ws.Range("E449").Value = 0.03347378277153558
$ws.Range("F449").Value = 0.0264516713470423
$ws.Range("G449").Value = 0.06413455318362443

$ws.Range("A450").ClearContents()
$ws.Range("B450").Value = "ROE"
$ws.Range("C450").Value = 0.03886763580719204
$ws.Range("D450").Value = 0.05066568047337278
$ws.Range("E450").Value = 0.03701307104956646
$ws.Range("F450").Value = 0.03903743315508022
$ws.Range("G450").Value = 0.09192052980132451

$ws.Range("A451").ClearContents()
$ws.Range("B451").Value = "ROA"
$ws.Range("C451").Value = 0.01457508463877891
$ws.Range("D451").Value = 0.01430112390827795
$ws.Range("E451").Value = 0.01038338658146965
$ws.Range("F451").Value = 0.009673358510567814
$ws.Range("G451").Value = 0.02769132551272843

